# Insert a new data row at row 197 (pushing the existing rows 197-264
# down to 198-265) and populate the new row with its values, matching
# the "Fruta / hortaliza, semanal" weekly-refresh edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 197..264 down to 198..265, leaving row 197 blank/ready.
$ws.Rows.Item(197).Insert()

# Populate the newly inserted row 197 with the new weekly record.
$ws.Cells.Item(197, 1).Value  = 8
$ws.Cells.Item(197, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(197, 3).Value  = "Coquimbo"
$ws.Cells.Item(197, 4).Value  = 45119
$ws.Cells.Item(197, 5).Value  = 4
$ws.Cells.Item(197, 6).Value  = 100112001
$ws.Cells.Item(197, 7).Value  = "Berenjena"
$ws.Cells.Item(197, 8).Value  = "Sin especificar"
$ws.Cells.Item(197, 9).Value  = "Primera"
$ws.Cells.Item(197, 10).Value = 320
$ws.Cells.Item(197, 11).Value = 8000
$ws.Cells.Item(197, 12).Value = 9000
$ws.Cells.Item(197, 13).Value = 8500
$ws.Cells.Item(197, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(197, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(197, 16).Value = 170
$ws.Cells.Item(197, 17).Value = 50
$ws.Cells.Item(197, 18).Value = "Hortaliza"
